# "Support basic group object"
# Adds a new slide 5 that demonstrates grouped shapes: a title, a group of
# three textboxes (the "top row" copied from slide 4) and a second group of
# three textboxes (the "bottom row" copied from slide 4), matching the
# target OOXML produced by PowerPoint's Group feature.

$p = $ppt.ActivePresentation
$s4 = $p.Slides.Item(4)

# --- Create the new slide (5th), using the same blank layout as slide 4 ---
# We briefly use a Title+Content layout so the "Content Placeholder" shape
# consumes id=3 (mirroring the id sequence of the original authoring
# session), then delete that placeholder, then swap in the real "blank"
# layout (slide4's layout) so the slide relationship matches slide 4's.
$s5 = $p.Slides.Add(5, 2)
$s5.Shapes.Item(2).Delete()
$s5.CustomLayout = $s4.CustomLayout

# --- Title ---
$titleSrc = $p.Slides.Item(2).Shapes.Item(1)
$titleSrc.Copy()
$titlePasted = $s5.Shapes.Paste()
$title = $titlePasted.Item(1)
$title.TextFrame.TextRange.Text = "群組測試"

# --- Top row: copy the three "top" textboxes from slide 4 ---
$s4.Shapes.Item(1).Copy()
$top1 = $s5.Shapes.Paste().Item(1)
$s4.Shapes.Item(2).Copy()
$top2 = $s5.Shapes.Paste().Item(1)
$s4.Shapes.Item(3).Copy()
$top3 = $s5.Shapes.Paste().Item(1)

$top1.Name = "文字方塊 3"
$top2.Name = "文字方塊 4"
$top3.Name = "文字方塊 5"

# Move the three shapes down (this is their position immediately before
# grouping -- PowerPoint freezes this as the group's child coordinate
# space, i.e. chOff/chExt).
$top1.Left = 72.50944881889764
$top1.Top  = 169.1016535433071
$top2.Left = 282.4371653543307
$top2.Top  = 169.1016535433071
$top3.Left = 492.36488188976375
$top3.Top  = 169.1016535433071

$topGroup = $s5.Shapes.Range(@($top1.Name, $top2.Name, $top3.Name)).Group()
$topGroup.Name = "群組 6"
# Nudge the whole group to its final resting place (off != chOff).
$topGroup.Left = 66.7951968503937
$topGroup.Top  = 169.1016535433071

# --- Bottom row: copy the three "bottom" textboxes from slide 4 ---
$s4.Shapes.Item(7).Copy()
$bot1 = $s5.Shapes.Paste().Item(1)
$s4.Shapes.Item(8).Copy()
$bot2 = $s5.Shapes.Paste().Item(1)
$s4.Shapes.Item(9).Copy()
$bot3 = $s5.Shapes.Paste().Item(1)

$bot1.Name = "文字方塊 7"
$bot2.Name = "文字方塊 8"
$bot3.Name = "文字方塊 9"

# These already land at slide4's original bottom-row coordinates, which is
# exactly the child coordinate space recorded by the diff, so no
# repositioning is required before grouping.

$bottomGroup = $s5.Shapes.Range(@($bot1.Name, $bot2.Name, $bot3.Name)).Group()
$bottomGroup.Name = "群組 10"
# Move the whole group to its final resting place.
$bottomGroup.Left = 0
$bottomGroup.Top  = 453.54330708661416
